$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Formula changes: switch the predicted-strain/stress block (columns AI, AP,
# AS, AU) from an incremental ("delta / running total") formulation to a
# total stress-strain formulation.
# ---------------------------------------------------------------------------

# AI: ratio used to build K_phi -> now derived from ev_pla/er_from_ev (AD/AE)
# instead of dev_pla/deq (AG/AH).
$ws.Range("AI3:AI18").Formula = "=-AD3/AE3"

# AP: slope of the swelling/compression line -> now a direct function of the
# already-computed ratio^4 value (AJ) instead of a secant slope between the
# current and previous rows.
$ws.Range("AP3:AP18").Formula = "=AJ3^0.25"

# AS: deviatoric plastic strain prediction -> now computed directly (total
# formulation, looking one row ahead) instead of accumulating a running sum
# of AR increments.
$ws.Range("AS2:AS18").Formula = "=(1+2*AM3)*AK3*(1-AP3/3)/(2*AN3*AO3*AM3)"

# AU: volumetric plastic strain prediction -> now a direct total-strain
# formula (based on the current row only) instead of a cumulative running
# total of AT increments.
$ws.Range("AU3:AU18").Formula = "=(1-AM3)*AK3*(1-AP3/3)/(3*AN3*AO3*AM3)"

$excel.Calculate()

# ---------------------------------------------------------------------------
# Highlight every cell touched by the reformulation in yellow, matching the
# workbook's existing "changed cell" convention.
# ---------------------------------------------------------------------------
$ws.Range("AI3:AI18").Interior.Color = 65535
$ws.Range("AP3:AP18").Interior.Color = 65535
$ws.Range("AS2:AS18").Interior.Color = 65535
$ws.Range("AU3:AU18").Interior.Color = 65535

# ---------------------------------------------------------------------------
# View state: re-centre/scroll the window and update zoom + selection to
# match where the author left off after the edit.
# ---------------------------------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 55
$win.ScrollColumn = 28
$win.ScrollRow = 1
$ws.Range("AO1:AO1048576").Select()

$excel.Calculate()
